$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first worksheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 288
$ws1.Range("F3").Value = 17
$ws1.Range("F4").Value = 7806
$ws1.Range("F5").Value = 5701
$ws1.Range("F7").Value = 78
$ws1.Range("F11").Value = 267

# Sheet "全部类型" (All types) - fourth worksheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 288
$ws4.Range("F3").Value = 17
$ws4.Range("F4").Value = 7806
$ws4.Range("F5").Value = 5701
$ws4.Range("F7").Value = 78
$ws4.Range("F13").Value = 267
